# Saldo_guide.xlsx edit script
# - Renames the worksheet tab to reflect the new export timestamp
# - Shifts every "Dt. Referencia" (column G) date by +1 day (45572 -> 45573)
# - Updates two rows whose "Saldo Previsto"/"Vl. Total" values changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet/tab (new export run name)
$ws.Name = "IClientBalance-20241008-090359-"

# Shift all reference dates in column G (rows 2-274) from 45572 to 45573
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value2 = 45573
}

# Row 108: Saldo Previsto / Vl. Total corrected from 70476.36 to 476.36
$ws.Cells.Item(108, 5).Value2 = 476.36
$ws.Cells.Item(108, 8).Value2 = 476.36

# Row 161: Saldo Previsto / Vl. Total corrected from 298.58999999999997 to 301.75
$ws.Cells.Item(161, 5).Value2 = 301.75
$ws.Cells.Item(161, 8).Value2 = 301.75
